$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked" (F) and "is_enabled" (G) header columns entirely.
# This shifts the "rem" column (formerly H) left into column F.
$ws.Range("F1:G1").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
